# Updated cryptos list - apply targeted cell value changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.616.89'
$ws.Range("D3").Value = '2.441.81'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '566.38'
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").Value = '145.81'
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '0.536'
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("E9").Value = '  +2.48%  '
$ws.Range("D10").Value = '0.155'
$ws.Range("E10").Value = '  +0.40%  '
$ws.Range("D11").Value = '5.26'
$ws.Range("E11").Value = '  -1.26%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '0.0000185'
$ws.Range("E13").Value = '  +5.64%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '26.87'
$ws.Range("E14").Value = '  +4.56%  '
$ws.Range("D15").Value = '2.834.04'
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("D16").Value = '62.425.83'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '2.438.42'
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").Value = '11.28'
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").Value = '6.96'
$ws.Range("E19").Value = '  +1.49%  '
$ws.Range("D20").Value = '324.75'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.40'
$ws.Range("E23").Value = '  +2.18%  '
$ws.Range("D24").Value = '1.74'
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("D25").Value = '8.77'
$ws.Range("E25").Value = '  -2.40%  '
$ws.Range("D26").Value = '0.0₃0999'
$ws.Range("E26").Value = '  +5.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '558.20'
$ws.Range("E27").Value = '  -3.68%  '
$ws.Range("D28").Value = '2.561.39'
$ws.Range("E28").Value = '  +1.20%  '
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D30").Value = '8.35'
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").Value = '1.46'
$ws.Range("E31").Value = '  +1.54%  '
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("D34").Value = '1.54'
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("D35").Value = '4.89'
$ws.Range("E35").Value = '  +3.36%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").Value = '0.383'
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("D39").Value = '18.85'
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("D40").Value = '150.28'
$ws.Range("E41").Value = '  +1.08%  '
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("D43").Value = '2.41'
$ws.Range("E43").Value = '  +5.30%  '
$ws.Range("D44").Value = '148.89'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '3.69'
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("E46").Value = '  +0.54%  '
$ws.Range("D47").Value = '20.52'
$ws.Range("E47").Value = '  +2.22%  '
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0930'
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '0.0232'
$ws.Range("E50").Value = '  +1.82%  '
$ws.Range("E51").Value = '  +0.40%  '
